# Convert the First-In / Last-Out Time columns (G,H) and the Date column (B)
# from numeric date/time serials into plain text strings (e.g. "1/1/2023",
# "5:00") for the three data rows, per the "time-based color code" commit.
#
# Excel's COM layer re-parses a string assigned to a cell that already has a
# date/time number format, turning it right back into a serial number. To
# force the literal text to stick we flip the cell's NumberFormat to the
# Text format ("@") before writing the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; Date = "1/1/2023"; FirstIn = "5:00";  LastOut = "13:00" },
    @{ Row = 3; Date = "1/2/2023"; FirstIn = "7:00";  LastOut = "15:00" },
    @{ Row = 4; Date = "1/3/2023"; FirstIn = "9:00";  LastOut = "17:00" }
)

foreach ($entry in $data) {
    $r = $entry.Row

    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $entry.Date

    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = $entry.FirstIn

    $ws.Cells.Item($r, 8).NumberFormat = "@"
    $ws.Cells.Item($r, 8).Value = $entry.LastOut
}
